$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L3").Value = 1.22
$ws.Range("M3").Value = 4
$ws.Range("N3").Value = 1.73
$ws.Range("O3").Value = 2.08
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 26
$ws.Range("AH4").Value = 126
$ws.Range("G4").Value = 1.25
$ws.Range("I4").Value = 9
$ws.Range("R4").Value = 1.91
$ws.Range("S4").Value = 1.91
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = 41
$ws.Range("AE6").Value = 11
$ws.Range("AF6").Value = 19
$ws.Range("AI6").Value = 29
$ws.Range("P6").Value = 1.36
$ws.Range("Q6").Value = 3
$ws.Range("Z6").Value = 11
$ws.Range("L7").Value = 1.25
$ws.Range("M7").Value = 3.75
$ws.Range("N7").Value = 1.9
$ws.Range("O7").Value = 1.95
$ws.Range("N8").Value = 1.83
$ws.Range("O8").Value = 2.03
$ws.Range("J9").Value = 1.02
$ws.Range("K9").Value = 19
$ws.Range("N9").Value = 1.5
$ws.Range("O9").Value = 2.5
$ws.Range("J10").Value = 1.02
$ws.Range("K10").Value = 21
$ws.Range("L10").Value = 1.11
$ws.Range("M10").Value = 6.5
$ws.Range("N10").Value = 1.4
$ws.Range("O10").Value = 2.88
$ws.Range("AB13").Value = 17
$ws.Range("AE13").Value = 7.3
$ws.Range("AF13").Value = 13.5
$ws.Range("Z13").Value = 7.1
$ws.Range("J19").Value = 1.04
$ws.Range("K19").Value = 13
$ws.Range("N19").Value = 1.85
$ws.Range("O19").Value = 2
$ws.Range("P19").Value = 1.36
$ws.Range("AE20").Value = 7
$ws.Range("G20").Value = 3.4
$ws.Range("I20").Value = 2.2
$ws.Range("L20").Value = 1.36
$ws.Range("M20").Value = 3
$ws.Range("P20").Value = 1.44
$ws.Range("Q20").Value = 2.63
$ws.Range("U20").Value = 17
$ws.Range("N21").Value = 1.93
$ws.Range("O21").Value = 1.93
$ws.Range("P21").Value = 1.36
$ws.Range("P22").Value = 1.5
$ws.Range("J23").Value = 1.06
$ws.Range("K23").Value = 10
$ws.Range("P23").Value = 1.44
$ws.Range("Q23").Value = 2.63
$ws.Range("G24").Value = 1.73
$ws.Range("H24").Value = 3.3
$ws.Range("I24").Value = 5.5
$ws.Range("K24").Value = 9.5
$ws.Range("P24").Value = 1.36
$ws.Range("Q24").Value = 3
$ws.Range("R24").Value = 1.8
$ws.Range("S24").Value = 1.95
$ws.Range("Y24").Value = 26
$ws.Range("Z24").Value = 9.5
$ws.Range("N26").Value = 1.9
$ws.Range("O26").Value = 1.95
$ws.Range("N27").Value = 2.2
$ws.Range("O27").Value = 1.65
$ws.Range("N28").Value = 1.7
$ws.Range("O28").Value = 2.1
$ws.Range("R29").Value = 2.16
$ws.Range("S29").Value = 1.62
$ws.Range("AA30").Value = 5.6
$ws.Range("AB30").Value = 14.5
$ws.Range("AC30").Value = 75
$ws.Range("AE30").Value = 8.25
$ws.Range("AF30").Value = 17.5
$ws.Range("AG30").Value = 12
$ws.Range("AH30").Value = 50
$ws.Range("AI30").Value = 35
$ws.Range("AJ30").Value = 45
$ws.Range("G30").Value = 1.78
$ws.Range("H30").Value = 3.25
$ws.Range("I30").Value = 4.2
$ws.Range("L30").Value = 1.5
$ws.Range("M30").Value = 2.27
$ws.Range("N30").Value = 2.12
$ws.Range("O30").Value = 1.57
$ws.Range("P30").Value = 1.42
$ws.Range("Q30").Value = 2.35
$ws.Range("R30").Value = 2.07
$ws.Range("S30").Value = 1.68
$ws.Range("T30").Value = 4.9
$ws.Range("U30").Value = 6.3
$ws.Range("V30").Value = 7.2
$ws.Range("W30").Value = 11.25
$ws.Range("X30").Value = 13
$ws.Range("Y30").Value = 27
$ws.Range("Z30").Value = 7.5
$ws.Range("AA33").Value = 13
$ws.Range("AB33").Value = 29
$ws.Range("AC33").Value = 81
$ws.Range("AE33").Value = 34
$ws.Range("AF33").Value = 67
$ws.Range("AG33").Value = 41
$ws.Range("AH33").Value = 201
$ws.Range("AI33").Value = 101
$ws.Range("AJ33").Value = 81
$ws.Range("G33").Value = 1.18
$ws.Range("H33").Value = 6.5
$ws.Range("J33").Value = 1.02
$ws.Range("K33").Value = 19
$ws.Range("R33").Value = 2.25
$ws.Range("S33").Value = 1.57
$ws.Range("T33").Value = 8.5
$ws.Range("U33").Value = 6
$ws.Range("V33").Value = 11
$ws.Range("Y33").Value = 34
